$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "45.294.74"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.426.22"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.01"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.50"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.77%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.55"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.21"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.08"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.807.44"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.417.40"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "45.227.01"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.33%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.98"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "244.79"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.28"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.62"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "49.52"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.00"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +7.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.01"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +7.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.21"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.67%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.44"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.47%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "126.66"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.08%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.65%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.73"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.16%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.936.99"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.01%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.76%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.79"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +7.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "76.55"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.81"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.38%  "
